$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.393.62'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -4.43%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.969.28'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -5.33%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '537.84'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -5.83%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.83'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -9.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.566'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.977.78'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -5.46%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -3.77%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.12'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -7.12%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -4.45%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.485.81'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -5.50%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.06%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '61.494.70'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -4.33%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '23.53'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -6.58%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.973.72'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -5.95%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -5.45%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.68%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.01'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -3.97%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '377.50'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -5.84%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -5.63%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -3.31%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '65.48'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -4.41%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.469'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -3.15%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.094.31'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -5.59%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.187'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -4.59%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.998'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.18%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0932'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -7.65%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.20'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -6.53%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -5.28%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '20.40'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -3.90%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '159.30'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.29%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.63'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -3.95%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.90'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -5.87%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.92%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -5.74%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.54'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -7.95%  '
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = 'OKB'
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '37.51'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.20%  '
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.90'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -3.85%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.402.75'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -8.99%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.08'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -6.70%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.668'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -3.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0589'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -4.03%  '
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'FirstDigitalUSD'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.997'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.03%  '
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.07'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -6.95%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0244'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -3.85%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0949'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.64%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.73'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -6.38%  '
